# Update EC (Estado de Cuenta) workbook:
#  - Remove two worker rows (JOHANA STEFANY MARTINEZ ANAYA, DEYNER ALFONSO LEAL PINEDA)
#    from the detail table, leaving ALFONSO JOSE SALGADO SILVA and MAURICIO ZUÑIGA TENA.
#  - Update the "Valor Mora" total and per-row amounts.
#  - Update worker/period counters.
#  - Re-fit column D (Nombre Trabajador) width to the shorter remaining names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the "DEYNER ALFONSO LEAL PINEDA" row (18) first, then the
# "JOHANA STEFANY MARTINEZ ANAYA" row (16), so row numbers above the
# second delete are unaffected by the first.
$ws.Rows("18:18").Delete()
$ws.Rows("16:16").Delete()

# Update the remaining detail rows' mora amounts.
$ws.Range("G16").Value = 3400000
$ws.Range("F16").Value = 46400
$ws.Range("G17").Value = 1160000
$ws.Range("F17").Value = 46400

# Update the summary header values.
$ws.Range("E11").Value = 92800
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1

# Column D no longer needs to fit the long names that were removed.
$ws.Columns("D:D").ColumnWidth = 28.95
